# Auto-generated edit script applying the scheduled-runner profit recalculations
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 135
$ws.Range("H135").Value = 1596
$ws.Range("I135").Value = 1427.8334
$ws.Range("K135").Value = 12850.5006
$ws.Range("M135").Value = -10315.5006

# Row 141
$ws.Range("H141").Value = 4488.1577
$ws.Range("I141").Value = 4316.353
$ws.Range("K141").Value = 12949.059
$ws.Range("M141").Value = -7769.059000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4505.5386
$ws.Range("I32").Value = 4466.2104
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 4466.2104
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -4179.2104
$ws.Range("N32").Value = -6574

# Row 55
$ws.Range("H55").Value = 28332.5
$ws.Range("J55").Value = 28332.5
$ws.Range("L55").Value = 28332.5
$ws.Range("N55").Value = -28962.5

# Row 61
$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -5424

# Row 124
$ws.Range("H124").Value = 10065.8
$ws.Range("J124").Value = 10065.8
$ws.Range("L124").Value = 10065.8
$ws.Range("N124").Value = -19885.8

# Row 132
$ws.Range("H132").Value = 1804.4
$ws.Range("I132").Value = 1005.5
$ws.Range("K132").Value = 3016.5
$ws.Range("M132").Value = -486.5

# Row 136
$ws.Range("H136").Value = 3500
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3191
$ws.Range("I134").Value = 3190.4375
$ws.Range("K134").Value = 9571.3125
$ws.Range("M134").Value = -7036.3125

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1796.2142
$ws.Range("I31").Value = 2000.1
$ws.Range("K31").Value = 2000.1
$ws.Range("M31").Value = -1705.1

# Row 34
$ws.Range("H34").Value = 1796.2142
$ws.Range("I34").Value = 2000.1
$ws.Range("K34").Value = 2000.1
$ws.Range("M34").Value = -1798.1

# Row 58
$ws.Range("H58").Value = 3277.6
$ws.Range("I58").Value = 2965
$ws.Range("J58").Value = 3746.5
$ws.Range("K58").Value = 2965
$ws.Range("L58").Value = 3746.5
$ws.Range("M58").Value = -2762
$ws.Range("N58").Value = -4152.5

# Row 94
$ws.Range("H94").Value = 2465.3333
$ws.Range("J94").Value = 2424
$ws.Range("L94").Value = 2424
$ws.Range("N94").Value = -3326

# Row 99
$ws.Range("H99").Value = 8900
$ws.Range("I99").Value = 8900
$ws.Range("K99").Value = 8900
$ws.Range("M99").Value = -7402

# Row 126
$ws.Range("H126").Value = 8900
$ws.Range("I126").Value = 8900
$ws.Range("K126").Value = 26700
$ws.Range("M126").Value = -24230

# Row 132
$ws.Range("H132").Value = 2199
$ws.Range("I132").Value = 2199
$ws.Range("K132").Value = 6597
$ws.Range("M132").Value = -4067

# Row 136
$ws.Range("H136").Value = 3277.6
$ws.Range("I136").Value = 2965
$ws.Range("J136").Value = 3746.5
$ws.Range("K136").Value = 8895
$ws.Range("L136").Value = 11239.5
$ws.Range("M136").Value = -6345
$ws.Range("N136").Value = -16339.5

# Row 141
$ws.Range("H141").Value = 96650
$ws.Range("J141").Value = 96650
$ws.Range("L141").Value = 96650
$ws.Range("N141").Value = -107010

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 73.27273
$ws.Range("I17").Value = 69
$ws.Range("J17").Value = 92.5
$ws.Range("K17").Value = 207
$ws.Range("L17").Value = 277.5
$ws.Range("M17").Value = -38
$ws.Range("N17").Value = -615.5

# Row 40
$ws.Range("H40").Value = 255.71428
$ws.Range("I40").Value = 138
$ws.Range("K40").Value = 552
$ws.Range("M40").Value = -483

# Row 92
$ws.Range("H92").Value = 899
$ws.Range("I92").Value = 900
$ws.Range("K92").Value = 2700
$ws.Range("M92").Value = -1452

# Row 129
$ws.Range("H129").Value = 2075.2727
$ws.Range("I129").Value = 747
$ws.Range("J129").Value = 3182.1667
$ws.Range("K129").Value = 2241
$ws.Range("L129").Value = 9546.500100000001
$ws.Range("M129").Value = 2759
$ws.Range("N129").Value = -19546.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 63
$ws.Range("H63").Value = 12333.333
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 65
$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 66
$ws.Range("H66").Value = 12333.333
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 70
$ws.Range("H70").Value = 4004
$ws.Range("I70").Value = 4004
$ws.Range("K70").Value = 4004
$ws.Range("M70").Value = -3734

# Row 73
$ws.Range("H73").Value = 4004
$ws.Range("I73").Value = 4004
$ws.Range("K73").Value = 4004
$ws.Range("M73").Value = -3068

# Row 126
$ws.Range("H126").Value = 1338.6666
$ws.Range("I126").Value = 1338.6666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4015.9998
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1545.9998
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Row 50
$ws.Range("H50").Value = 43042
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 43042
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 43042
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -44316

# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()

# Row 61
$ws.Range("H61").Value = 1666.5
$ws.Range("I61").Value = 1666.5
$ws.Range("K61").Value = 1666.5
$ws.Range("M61").Value = -1464.5

# Row 63
$ws.Range("H63").Value = 29996.666
$ws.Range("I63").Value = 29990
$ws.Range("K63").Value = 29990
$ws.Range("M63").Value = -29241

# Row 66
$ws.Range("H66").Value = 29996.666
$ws.Range("I66").Value = 29990
$ws.Range("K66").Value = 89970
$ws.Range("M66").Value = -86226

# Row 82
$ws.Range("H82").Value = 1525.2727
$ws.Range("I82").Value = 1235.6
$ws.Range("K82").Value = 1235.6
$ws.Range("M82").Value = -874.5999999999999

# Row 85
$ws.Range("H85").Value = 1525.2727
$ws.Range("I85").Value = 1235.6
$ws.Range("K85").Value = 1235.6
$ws.Range("M85").Value = 12.40000000000009

# Row 113
$ws.Range("H113").Value = 1666.5
$ws.Range("I113").Value = 1666.5
$ws.Range("K113").Value = 1666.5
$ws.Range("M113").Value = 503.5

# Row 127
$ws.Range("H127").Value = 29107.5
$ws.Range("J127").Value = 29107.5
$ws.Range("L127").Value = 29107.5
$ws.Range("N127").Value = -39027.5

$ws = $wb.Worksheets.Item("WVR")
# Row 38
$ws.Range("H38").Value = 50000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 50000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 50000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -50946

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 51
$ws.Range("H51").Value = 34499.5
$ws.Range("I51").Value = 28999
$ws.Range("K51").Value = 28999
$ws.Range("M51").Value = -28489

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 58
$ws.Range("H58").Value = 18783.334
$ws.Range("I58").Value = 5675
$ws.Range("K58").Value = 5675
$ws.Range("M58").Value = -5367

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 107
$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -6840

# Row 126
$ws.Range("H126").Value = 3757.9092
$ws.Range("I126").Value = 3354.625
$ws.Range("K126").Value = 10063.875
$ws.Range("M126").Value = -7593.875

# Row 136
$ws.Range("H136").Value = 1617.3846
$ws.Range("I136").Value = 1617.3846
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4852.1538
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2302.1538
$ws.Range("N136").ClearContents()
